$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("H2").Value = 2
$ws1.Range("L2").Value = 1.12

# Row 3
$ws1.Range("H3").Value = 1
$ws1.Range("L3").Value = 1.1

# Row 4
$ws1.Range("H4").Value = 0
$ws1.Range("I4").Value = "High"
$ws1.Range("J4").Value = "Urgent"
$ws1.Range("L4").Value = 0.88

# Row 5
$ws1.Range("H5").Value = 0
$ws1.Range("L5").Value = 1.15

# Row 6
$ws1.Range("L6").Value = 1.12

# Row 7
$ws1.Range("L7").Value = 1.01

# Row 8
$ws1.Range("L8").Value = 1.05

# Row 9
$ws1.Range("L9").Value = 1.01

# Row 10
$ws1.Range("L10").Value = 1.17

# Row 11
$ws1.Range("L11").Value = 0.9

# Row 12
$ws1.Range("L12").Value = 1.2

# Row 13
$ws1.Range("L13").Value = 1.17

# Row 14
$ws1.Range("L14").Value = 0.97

# Row 15
$ws1.Range("L15").Value = 0.89

# Row 16
$ws1.Range("L16").Value = 1.14

# Row 17
$ws1.Range("L17").Value = 1.08

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "10"
$ws2.Range("B10").Value = "4"
$ws2.Range("B11").Value = "2"
